$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "26.297.88"
$ws.Range("E2").Value2 = "  +0.79%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "1.664.13"
$ws.Range("E3").Value2 = "  +0.55%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "219.05"
$ws.Range("E5").Value2 = "  +0.68%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "0.5339"
$ws.Range("E6").Value2 = "  +1.42%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.2656"
$ws.Range("E8").Value2 = "  +1.80%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.06416"
$ws.Range("E9").Value2 = "  +1.16%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "20.60"
$ws.Range("E10").Value2 = "  +0.92%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.07841"
$ws.Range("E11").Value2 = "  +0.62%  "

# Row 12
$ws.Range("E12").Value2 = "  +1.55%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "1.665.30"
$ws.Range("E13").Value2 = "  +4.42%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "1.892.22"
$ws.Range("E14").Value2 = "  +0.45%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "0.5529"
$ws.Range("E15").Value2 = "  +0.72%  "

# Row 16
$ws.Range("E16").Value2 = "  -0.04%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "65.72"
$ws.Range("E17").Value2 = "  +0.50%  "

# Row 18
$ws.Range("B18").Value2 = "Dai"
$ws.Range("C18").Value2 = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "1.010"
$ws.Range("E18").Value2 = "  +0.84%  "

# Row 19
$ws.Range("B19").Value2 = "Uniswap"
$ws.Range("C19").Value2 = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "4.690"
$ws.Range("E19").Value2 = "  +2.46%  "

# Row 20
$ws.Range("B20").Value2 = "BitcoinCash"
$ws.Range("C20").Value2 = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "194.33"
$ws.Range("E20").Value2 = "  +1.94%  "

# Row 21
$ws.Range("B21").Value2 = "Avalanche"
$ws.Range("C21").Value2 = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "10.22"
$ws.Range("E21").Value2 = "  +1.55%  "

# Row 22
$ws.Range("B22").Value2 = "Chainlink"
$ws.Range("C22").Value2 = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "6.045"
$ws.Range("E22").Value2 = "  +0.22%  "

# Row 23
$ws.Range("B23").Value2 = "BinanceUSD"
$ws.Range("C23").Value2 = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "1.011"
$ws.Range("E23").Value2 = "  +0.80%  "

# Row 24
$ws.Range("B24").Value2 = "Monero"
$ws.Range("C24").Value2 = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "146.13"
$ws.Range("E24").Value2 = "  +3.19%  "

# Row 25
$ws.Range("B25").Value2 = "Stellar"
$ws.Range("C25").Value2 = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "0.1235"
$ws.Range("E25").Value2 = "  -0.05%  "

# Row 26
$ws.Range("B26").Value2 = "Cosmos"
$ws.Range("C26").Value2 = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "7.202"
$ws.Range("E26").Value2 = "  -0.52%  "

# Row 27
$ws.Range("B27").Value2 = "EthereumClassic"
$ws.Range("C27").Value2 = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "16.18"
$ws.Range("E27").Value2 = "  +0.77%  "

# Row 28
$ws.Range("B28").Value2 = "Toncoin"
$ws.Range("C28").Value2 = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "1.483"
$ws.Range("E28").Value2 = "  +3.88%  "

# Row 29
$ws.Range("B29").Value2 = "Hedera"
$ws.Range("C29").Value2 = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "0.05855"
$ws.Range("E29").Value2 = "  -0.48%  "

# Row 30
$ws.Range("B30").Value2 = "PancakeSwap"
$ws.Range("C30").Value2 = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "1.281"
$ws.Range("E30").Value2 = "  +0.50%  "

# Row 31
$ws.Range("B31").Value2 = "InternetComputer(DFINITY)"
$ws.Range("C31").Value2 = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "3.625"
$ws.Range("E31").Value2 = "  +2.97%  "

# Row 32
$ws.Range("B32").Value2 = "Filecoin"
$ws.Range("C32").Value2 = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "3.284"
$ws.Range("E32").Value2 = "  +0.69%  "

# Row 33
$ws.Range("B33").Value2 = "LidoDAOToken"
$ws.Range("C33").Value2 = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "1.613"
$ws.Range("E33").Value2 = "  +1.42%  "

# Row 34
$ws.Range("B34").Value2 = "ARBITRUM"
$ws.Range("C34").Value2 = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "0.9645"
$ws.Range("E34").Value2 = "  +1.40%  "

# Row 35
$ws.Range("B35").Value2 = "MXToken"
$ws.Range("C35").Value2 = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "2.828"
$ws.Range("E35").Value2 = "  +1.57%  "

# Row 36
$ws.Range("B36").Value2 = "HuobiToken"
$ws.Range("C36").Value2 = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "2.417"
$ws.Range("E36").Value2 = "  +0.21%  "

# Row 37
$ws.Range("B37").Value2 = "ImmutableX"
$ws.Range("C37").Value2 = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "0.5809"
$ws.Range("E37").Value2 = "  +2.02%  "

# Row 38
$ws.Range("B38").Value2 = "VeChain"
$ws.Range("C38").Value2 = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "0.01609"
$ws.Range("E38").Value2 = "  -0.34%  "

# Row 39
$ws.Range("B39").Value2 = "TrustWalletToken"
$ws.Range("C39").Value2 = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "0.8674"
$ws.Range("E39").Value2 = "  +2.27%  "

# Row 40
$ws.Range("B40").Value2 = "FraxShare"
$ws.Range("C40").Value2 = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "5.886"
$ws.Range("E40").Value2 = "  +1.52%  "

# Row 41
$ws.Range("B41").Value2 = "Maker"
$ws.Range("C41").Value2 = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "1.052.51"
$ws.Range("E41").Value2 = "  +2.52%  "

# Row 42
$ws.Range("B42").Value2 = "PaxDollar"
$ws.Range("C42").Value2 = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "1.010"
$ws.Range("E42").Value2 = "  +0.78%  "

# Row 43
$ws.Range("B43").Value2 = "Quant"
$ws.Range("C43").Value2 = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "104.57"
$ws.Range("E43").Value2 = "  +1.98%  "

# Row 44
$ws.Range("B44").Value2 = "RocketPoolETH"
$ws.Range("C44").Value2 = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "1.802.84"
$ws.Range("E44").Value2 = "  +0.27%  "

# Row 45
$ws.Range("B45").Value2 = "Aave"
$ws.Range("C45").Value2 = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "57.77"
$ws.Range("E45").Value2 = "  +1.20%  "

# Row 46
$ws.Range("B46").Value2 = "BabyDogeCoin"
$ws.Range("C46").Value2 = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "0.0$([char]0x2088)106"
$ws.Range("E46").Value2 = "  -5.33%  "

# Row 47
$ws.Range("B47").Value2 = "Frax"
$ws.Range("C47").Value2 = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "1.012"
$ws.Range("E47").Value2 = "  +1.13%  "

# Row 48
$ws.Range("B48").Value2 = "Mantle"
$ws.Range("C48").Value2 = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "0.4384"
$ws.Range("E48").Value2 = "  +1.99%  "

# Row 49
$ws.Range("B49").Value2 = "EnergySwap"
$ws.Range("C49").Value2 = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "8.039"
$ws.Range("E49").Value2 = "  +2.42%  "

# Row 50
$ws.Range("B50").Value2 = "Cronos"
$ws.Range("C50").Value2 = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "0.05165"
$ws.Range("E50").Value2 = "  +0.25%  "

# Row 51
$ws.Range("B51").Value2 = "RenderToken"
$ws.Range("C51").Value2 = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "1.418"
$ws.Range("E51").Value2 = "  -3.87%  "
